$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 185, shifting existing rows 185-286 down to 186-287
$ws.Rows("185:185").Insert()

# Populate the newly inserted row 185 with the new weekly price record
$ws.Range("A185").Value = 9
$ws.Range("B185").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C185").Value = "Metropolitana"
$ws.Range("D185").Value = 44572
$ws.Range("E185").Value = 13
$ws.Range("F185").Value = 100112044
$ws.Range("G185").Value = "Perejil"
$ws.Range("H185").Value = "Sin especificar"
$ws.Range("I185").Value = "Primera"
$ws.Range("J185").Value = 61
$ws.Range("K185").Value = 14000
$ws.Range("L185").Value = 16000
$ws.Range("M185").Value = 15016
$ws.Range("N185").Value = "$/docena de atados"
$ws.Range("O185").Value = "Región Metropolitana"
$ws.Range("P185").Value = 5005
$ws.Range("Q185").Value = 3
$ws.Range("R185").Value = "Hortaliza"
